{"js": "// Update the two-digit multiplication expressions in the table to the\n// newly generated set of problems. Each old expression is unique in the\n// document, so a simple search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"94\u00d772=\", \"75\u00d784=\"],\n  [\"51\u00d749=\", \"76\u00d737=\"],\n  [\"49\u00d742=\", \"18\u00d751=\"],\n  [\"96\u00d760=\", \"52\u00d761=\"],\n  [\"78\u00d752=\", \"57\u00d770=\"],\n  [\"11\u00d734=\", \"97\u00d763=\"],\n  [\"25\u00d783=\", \"28\u00d751=\"],\n  [\"17\u00d766=\", \"50\u00d790=\"],\n  [\"19\u00d770=\", \"62\u00d758=\"],\n  [\"63\u00d765=\", \"46\u00d799=\"],\n  [\"69\u00d743=\", \"18\u00d780=\"],\n  [\"47\u00d797=\", \"59\u00d799=\"],\n  [\"32\u00d717=\", \"25\u00d742=\"],\n  [\"32\u00d746=\", \"68\u00d795=\"],\n  [\"75\u00d781=\", \"17\u00d716=\"],\n  [\"46\u00d714=\", \"29\u00d745=\"],\n  [\"90\u00d755=\", \"39\u00d774=\"],\n  [\"26\u00d750=\", \"35\u00d789=\"],\n  [\"47\u00d744=\", \"21\u00d734=\"],\n  [\"90\u00d735=\", \"54\u00d775=\"],\n  [\"66\u00d748=\", \"52\u00d769=\"],\n  [\"19\u00d755=\", \"25\u00d732=\"],\n  [\"24\u00d773=\", \"73\u00d767=\"],\n  [\"57\u00d758=\", \"30\u00d739=\"],\n  [\"19\u00d775=\", \"94\u00d748=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication expressions in the table to the\n# newly generated set of problems. Each old expression is unique in the\n# document, so Find/Execute per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"94\u00d772=\", \"75\u00d784=\"),\n    @(\"51\u00d749=\", \"76\u00d737=\"),\n    @(\"49\u00d742=\", \"18\u00d751=\"),\n    @(\"96\u00d760=\", \"52\u00d761=\"),\n    @(\"78\u00d752=\", \"57\u00d770=\"),\n    @(\"11\u00d734=\", \"97\u00d763=\"),\n    @(\"25\u00d783=\", \"28\u00d751=\"),\n    @(\"17\u00d766=\", \"50\u00d790=\"),\n    @(\"19\u00d770=\", \"62\u00d758=\"),\n    @(\"63\u00d765=\", \"46\u00d799=\"),\n    @(\"69\u00d743=\", \"18\u00d780=\"),\n    @(\"47\u00d797=\", \"59\u00d799=\"),\n    @(\"32\u00d717=\", \"25\u00d742=\"),\n    @(\"32\u00d746=\", \"68\u00d795=\"),\n    @(\"75\u00d781=\", \"17\u00d716=\"),\n    @(\"46\u00d714=\", \"29\u00d745=\"),\n    @(\"90\u00d755=\", \"39\u00d774=\"),\n    @(\"26\u00d750=\", \"35\u00d789=\"),\n    @(\"47\u00d744=\", \"21\u00d734=\"),\n    @(\"90\u00d735=\", \"54\u00d775=\"),\n    @(\"66\u00d748=\", \"52\u00d769=\"),\n    @(\"19\u00d755=\", \"25\u00d732=\"),\n    @(\"24\u00d773=\", \"73\u00d767=\"),\n    @(\"57\u00d758=\", \"30\u00d739=\"),\n    @(\"19\u00d775=\", \"94\u00d748=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
